# Add two new Mac-Address rows to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: regcntr_id=10001, usr_id=110030, lang_code=eng, is_active=TRUE, cr_by=superadmin, cr_dtimes=now(), eff_dtimes=now()
$ws.Cells.Item(31, 1).Value = 10001
$ws.Cells.Item(31, 2).Value = 110030
$ws.Cells.Item(31, 3).Value = "eng"
$ws.Cells.Item(31, 4).Value = $true
$ws.Cells.Item(31, 5).Value = "superadmin"
$ws.Cells.Item(31, 6).Value = "now()"
$ws.Cells.Item(31, 7).Value = "now()"

# Row 32: regcntr_id=10001, usr_id=110031, lang_code=eng, is_active=TRUE, cr_by=superadmin, cr_dtimes=now(), eff_dtimes=now()
$ws.Cells.Item(32, 1).Value = 10001
$ws.Cells.Item(32, 2).Value = 110031
$ws.Cells.Item(32, 3).Value = "eng"
$ws.Cells.Item(32, 4).Value = $true
$ws.Cells.Item(32, 5).Value = "superadmin"
$ws.Cells.Item(32, 6).Value = "now()"
$ws.Cells.Item(32, 7).Value = "now()"

# Update the view selection/top-left cell to reflect the new active area
$ws.Range("E28").Select()
$excel.ActiveWindow.ScrollRow = 19
